# ---------------------------------------------------------------------------
# Reporting-format improvements: split the "QUERY" block into two blocks
# ("QUERY - Corrected" in R:V, "QUERY - Original" in the newly added W:AA),
# and add a "TOTAL REAL TIME" series to the two line charts that plot
# per-operation-size timings (QUERY chart + SEARCH chart).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Duplicate the existing QUERY columns (R:V) into the new columns (W:AA),
#    preserving the original ("uncorrected") numbers, and relabel the two
#    header blocks.
# ---------------------------------------------------------------------------

# Copy format + values of the R:V block (header + data rows) into W:AA first,
# so the new columns inherit the same look (bold/centred header, plain data).
$ws.Range("R1:V9").Copy()
$ws.Range("W1").PasteSpecial(-4104)

# W1:AA1 needs to be its own merged header cell (like R1:V1 already is).
$ws.Range("W1:AA1").Merge()

# Re-label the two header blocks.
$ws.Range("R1").Value() = "QUERY - Corrected"
$ws.Range("W1").Value() = "QUERY - Original"

# ---------------------------------------------------------------------------
# 2. Apply the "corrected" numbers to the QUERY - Corrected block (only the
#    1-thread and 10-thread TOTAL CPU TIME figures actually changed; the
#    other rows/columns stay as they were).
# ---------------------------------------------------------------------------
$ws.Range("R3").Value() = 0
$ws.Range("T3").Value() = 0

# ---------------------------------------------------------------------------
# 3. Resize / reposition the QUERY and SEARCH line charts and add a
#    "TOTAL REAL TIME" series (sourced from row 3) to each, inserted as the
#    first series in the plot order.
# ---------------------------------------------------------------------------

# --- QUERY chart (now plots Sheet1!$R block => "QUERY - Corrected") -------
$queryChart = $ws.ChartObjects().Item(2)
$queryChart.Left = 1580.6748046875
$queryChart.Top = 339.0
$queryChart.Width = 485.4238281249997
$queryChart.Height = 422.0

$querySeries = $queryChart.Chart.SeriesCollection().NewSeries()
$querySeries.Name = "=Sheet1!`$G`$3"
$querySeries.Values = "=Sheet1!`$R`$3:`$V`$3"
$querySeries.XValues = "=Sheet1!`$R`$2:`$V`$2"
$querySeries.PlotOrder = 1

# --- SEARCH chart (plots Sheet1!$M block) ----------------------------------
$searchChart = $ws.ChartObjects().Item(3)
$searchChart.Left = 1075.1748046875
$searchChart.Top = 337.5
$searchChart.Width = 487.0
$searchChart.Height = 424.5

$searchSeries = $searchChart.Chart.SeriesCollection().NewSeries()
$searchSeries.Name = "=Sheet1!`$G`$3"
$searchSeries.Values = "=Sheet1!`$M`$3:`$Q`$3"
$searchSeries.XValues = "=Sheet1!`$M`$2:`$Q`$2"
$searchSeries.PlotOrder = 1

# ---------------------------------------------------------------------------
# 4. View cosmetics matching the refreshed layout.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.Zoom = 155
